# Atualizacao de bases das ligas, do dia: 23-02-2024 as 23:34
# Fixes a row-pairing mixup: several adjacent match rows had their
# betting-odds data swapped; this restores each row (id/teams/odds)
# to its correct match while keeping the row's own sequential index (col A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("B17").Value2 = 6221766
$ws.Range("F17").Value2 = 'Kairat Almaty'
$ws.Range("G17").Value2 = 'FK Kaspyi Aktau'
$ws.Range("H17").Value2 = 3
$ws.Range("I17").Value2 = 1
$ws.Range("J17").Value2 = 'H'
$ws.Range("K17").Value2 = 1.55
$ws.Range("L17").Value2 = 3.8
$ws.Range("M17").Value2 = 5
$ws.Range("N17").Value2 = 1.65
$ws.Range("O17").Value2 = 4
$ws.Range("P17").Value2 = 4.5
$ws.Range("Q17").Value2 = -0.75
$ws.Range("R17").Value2 = 1.8
$ws.Range("S17").Value2 = 2
$ws.Range("T17").Value2 = 2.75
$ws.Range("U17").Value2 = 1.925
$ws.Range("V17").Value2 = 1.875
$ws.Range("W17").Value2 = 0.6499999999999999
$ws.Range("X17").Value2 = -1
$ws.Range("Y17").Value2 = -1
$ws.Range("Z17").Value2 = 0.8
$ws.Range("AA17").Value2 = -1
$ws.Range("AB17").Value2 = 0.925
$ws.Range("AC17").Value2 = -1

# Row 18
$ws.Range("B18").Value2 = 6221641
$ws.Range("F18").Value2 = 'Tobol Kostanay'
$ws.Range("G18").Value2 = 'Shakhter Karagandy'
$ws.Range("H18").Value2 = 2
$ws.Range("I18").Value2 = 1
$ws.Range("J18").Value2 = 'H'
$ws.Range("K18").Value2 = 1.4
$ws.Range("L18").Value2 = 4.333
$ws.Range("M18").Value2 = 6
$ws.Range("N18").Value2 = 1.333
$ws.Range("O18").Value2 = 4.75
$ws.Range("P18").Value2 = 6.5
$ws.Range("Q18").Value2 = -1.5
$ws.Range("R18").Value2 = 1.925
$ws.Range("S18").Value2 = 1.875
$ws.Range("T18").Value2 = 2.75
$ws.Range("U18").Value2 = 1.75
$ws.Range("V18").Value2 = 1.95
$ws.Range("W18").Value2 = 0.333
$ws.Range("X18").Value2 = -1
$ws.Range("Y18").Value2 = -1
$ws.Range("Z18").Value2 = -1
$ws.Range("AA18").Value2 = 0.875
$ws.Range("AB18").Value2 = 0.375
$ws.Range("AC18").Value2 = -0.5

# Row 47
$ws.Range("B47").Value2 = 6221661
$ws.Range("F47").Value2 = 'FK Atyrau'
$ws.Range("G47").Value2 = 'FK Aksu'
$ws.Range("H47").Value2 = 2
$ws.Range("I47").Value2 = 1
$ws.Range("J47").Value2 = 'H'
$ws.Range("K47").Value2 = 2
$ws.Range("L47").Value2 = 3.4
$ws.Range("M47").Value2 = 3.1
$ws.Range("N47").Value2 = 1.833
$ws.Range("O47").Value2 = 3.5
$ws.Range("P47").Value2 = 3.5
$ws.Range("Q47").Value2 = -0.5
$ws.Range("R47").Value2 = 1.9
$ws.Range("S47").Value2 = 1.9
$ws.Range("T47").Value2 = 2.25
$ws.Range("U47").Value2 = 1.95
$ws.Range("V47").Value2 = 1.85
$ws.Range("W47").Value2 = 0.833
$ws.Range("X47").Value2 = -1
$ws.Range("Y47").Value2 = -1
$ws.Range("Z47").Value2 = 0.8999999999999999
$ws.Range("AA47").Value2 = -1
$ws.Range("AB47").Value2 = 0.95
$ws.Range("AC47").Value2 = -1

# Row 48
$ws.Range("B48").Value2 = 6221658
$ws.Range("F48").Value2 = 'Tobol Kostanay'
$ws.Range("G48").Value2 = 'FK Aktobe'
$ws.Range("H48").Value2 = 2
$ws.Range("I48").Value2 = 2
$ws.Range("J48").Value2 = 'D'
$ws.Range("K48").Value2 = 2.4
$ws.Range("L48").Value2 = 3.2
$ws.Range("M48").Value2 = 2.625
$ws.Range("N48").Value2 = 2.4
$ws.Range("O48").Value2 = 3.1
$ws.Range("P48").Value2 = 2.7
$ws.Range("Q48").Value2 = 0
$ws.Range("R48").Value2 = 1.775
$ws.Range("S48").Value2 = 2.025
$ws.Range("T48").Value2 = 2.25
$ws.Range("U48").Value2 = 1.775
$ws.Range("V48").Value2 = 2.025
$ws.Range("W48").Value2 = -1
$ws.Range("X48").Value2 = 2.1
$ws.Range("Y48").Value2 = -1
$ws.Range("Z48").Value2 = 0
$ws.Range("AA48").Value2 = -0
$ws.Range("AB48").Value2 = 0.7749999999999999
$ws.Range("AC48").Value2 = -1

# Row 51
$ws.Range("B51").Value2 = 6221666
$ws.Range("F51").Value2 = 'Tobol Kostanay'
$ws.Range("G51").Value2 = 'Kairat Almaty'
$ws.Range("H51").Value2 = 2
$ws.Range("I51").Value2 = 3
$ws.Range("J51").Value2 = 'A'
$ws.Range("K51").Value2 = 2.05
$ws.Range("L51").Value2 = 3.4
$ws.Range("M51").Value2 = 3.2
$ws.Range("N51").Value2 = 2.1
$ws.Range("O51").Value2 = 3.4
$ws.Range("P51").Value2 = 3
$ws.Range("Q51").Value2 = -0.25
$ws.Range("R51").Value2 = 1.875
$ws.Range("S51").Value2 = 1.925
$ws.Range("T51").Value2 = 2.5
$ws.Range("U51").Value2 = 1.8
$ws.Range("V51").Value2 = 2
$ws.Range("W51").Value2 = -1
$ws.Range("X51").Value2 = -1
$ws.Range("Y51").Value2 = 2
$ws.Range("Z51").Value2 = -1
$ws.Range("AA51").Value2 = 0.925
$ws.Range("AB51").Value2 = 0.8
$ws.Range("AC51").Value2 = -1

# Row 52
$ws.Range("B52").Value2 = 6221662
$ws.Range("F52").Value2 = 'FK Aktobe'
$ws.Range("G52").Value2 = 'FK Aksu'
$ws.Range("H52").Value2 = 2
$ws.Range("I52").Value2 = 1
$ws.Range("J52").Value2 = 'H'
$ws.Range("K52").Value2 = 1.2
$ws.Range("L52").Value2 = 6
$ws.Range("M52").Value2 = 10
$ws.Range("N52").Value2 = 1.285
$ws.Range("O52").Value2 = 6
$ws.Range("P52").Value2 = 6.5
$ws.Range("Q52").Value2 = -1.75
$ws.Range("R52").Value2 = 1.925
$ws.Range("S52").Value2 = 1.875
$ws.Range("T52").Value2 = 3.25
$ws.Range("U52").Value2 = 2
$ws.Range("V52").Value2 = 1.8
$ws.Range("W52").Value2 = 0.2849999999999999
$ws.Range("X52").Value2 = -1
$ws.Range("Y52").Value2 = -1
$ws.Range("Z52").Value2 = -1
$ws.Range("AA52").Value2 = 0.875
$ws.Range("AB52").Value2 = -0.5
$ws.Range("AC52").Value2 = 0.4

# Row 66
$ws.Range("B66").Value2 = 6221673
$ws.Range("F66").Value2 = 'Shakhter Karagandy'
$ws.Range("G66").Value2 = 'FK Aksu'
$ws.Range("H66").Value2 = 2
$ws.Range("I66").Value2 = 1
$ws.Range("J66").Value2 = 'H'
$ws.Range("K66").Value2 = 2.1
$ws.Range("L66").Value2 = 3.25
$ws.Range("M66").Value2 = 3
$ws.Range("N66").Value2 = 2.1
$ws.Range("O66").Value2 = 3.3
$ws.Range("P66").Value2 = 3
$ws.Range("Q66").Value2 = -0.25
$ws.Range("R66").Value2 = 1.925
$ws.Range("S66").Value2 = 1.875
$ws.Range("T66").Value2 = 2.5
$ws.Range("U66").Value2 = 1.975
$ws.Range("V66").Value2 = 1.825
$ws.Range("W66").Value2 = 1.1
$ws.Range("X66").Value2 = -1
$ws.Range("Y66").Value2 = -1
$ws.Range("Z66").Value2 = 0.925
$ws.Range("AA66").Value2 = -1
$ws.Range("AB66").Value2 = 0.9750000000000001
$ws.Range("AC66").Value2 = -1

# Row 67
$ws.Range("B67").Value2 = 6221674
$ws.Range("F67").Value2 = 'Zhetysu'
$ws.Range("G67").Value2 = 'FK Atyrau'
$ws.Range("H67").Value2 = 0
$ws.Range("I67").Value2 = 3
$ws.Range("J67").Value2 = 'A'
$ws.Range("K67").Value2 = 2.5
$ws.Range("L67").Value2 = 3.2
$ws.Range("M67").Value2 = 2.5
$ws.Range("N67").Value2 = 2.5
$ws.Range("O67").Value2 = 3.2
$ws.Range("P67").Value2 = 2.5
$ws.Range("Q67").Value2 = 0
$ws.Range("R67").Value2 = 1.9
$ws.Range("S67").Value2 = 1.9
$ws.Range("T67").Value2 = 2.25
$ws.Range("U67").Value2 = 1.925
$ws.Range("V67").Value2 = 1.875
$ws.Range("W67").Value2 = -1
$ws.Range("X67").Value2 = -1
$ws.Range("Y67").Value2 = 1.5
$ws.Range("Z67").Value2 = -1
$ws.Range("AA67").Value2 = 0.8999999999999999
$ws.Range("AB67").Value2 = 0.925
$ws.Range("AC67").Value2 = -1

# Row 82
$ws.Range("B82").Value2 = 6726053
$ws.Range("F82").Value2 = 'Tobol Kostanay'
$ws.Range("G82").Value2 = 'FK Kyzylzhar'
$ws.Range("H82").Value2 = 1
$ws.Range("I82").Value2 = 0
$ws.Range("J82").Value2 = 'H'
$ws.Range("K82").Value2 = 2.5
$ws.Range("L82").Value2 = 3.2
$ws.Range("M82").Value2 = 2.5
$ws.Range("N82").Value2 = 1.8
$ws.Range("O82").Value2 = 3
$ws.Range("P82").Value2 = 4.2
$ws.Range("Q82").Value2 = -0.5
$ws.Range("R82").Value2 = 1.9
$ws.Range("S82").Value2 = 1.9
$ws.Range("T82").Value2 = 2.25
$ws.Range("U82").Value2 = 1.875
$ws.Range("V82").Value2 = 1.925
$ws.Range("W82").Value2 = 0.8
$ws.Range("X82").Value2 = -1
$ws.Range("Y82").Value2 = -1
$ws.Range("Z82").Value2 = 0.8999999999999999
$ws.Range("AA82").Value2 = -1
$ws.Range("AB82").Value2 = -1
$ws.Range("AC82").Value2 = 0.925

# Row 83
$ws.Range("B83").Value2 = 6221684
$ws.Range("F83").Value2 = 'FC Astana'
$ws.Range("G83").Value2 = 'FK Aksu'
$ws.Range("H83").Value2 = 1
$ws.Range("I83").Value2 = 0
$ws.Range("J83").Value2 = 'H'
$ws.Range("K83").Value2 = 1.533
$ws.Range("L83").Value2 = 4.2
$ws.Range("M83").Value2 = 4.5
$ws.Range("N83").Value2 = 1.222
$ws.Range("O83").Value2 = 5.5
$ws.Range("P83").Value2 = 9
$ws.Range("Q83").Value2 = -2
$ws.Range("R83").Value2 = 1.975
$ws.Range("S83").Value2 = 1.825
$ws.Range("T83").Value2 = 3.25
$ws.Range("U83").Value2 = 1.925
$ws.Range("V83").Value2 = 1.875
$ws.Range("W83").Value2 = 0.222
$ws.Range("X83").Value2 = -1
$ws.Range("Y83").Value2 = -1
$ws.Range("Z83").Value2 = -1
$ws.Range("AA83").Value2 = 0.825
$ws.Range("AB83").Value2 = -1
$ws.Range("AC83").Value2 = 0.875

# Row 95
$ws.Range("B95").Value2 = 6221698
$ws.Range("F95").Value2 = 'FK Maktaaral'
$ws.Range("G95").Value2 = 'FK Aktobe'
$ws.Range("H95").Value2 = 1
$ws.Range("I95").Value2 = 2
$ws.Range("J95").Value2 = 'A'
$ws.Range("K95").Value2 = 4.333
$ws.Range("L95").Value2 = 3.5
$ws.Range("M95").Value2 = 1.666
$ws.Range("N95").Value2 = 4.2
$ws.Range("O95").Value2 = 3.4
$ws.Range("P95").Value2 = 1.7
$ws.Range("Q95").Value2 = 0.75
$ws.Range("R95").Value2 = 1.825
$ws.Range("S95").Value2 = 1.975
$ws.Range("T95").Value2 = 2.5
$ws.Range("U95").Value2 = 1.925
$ws.Range("V95").Value2 = 1.875
$ws.Range("W95").Value2 = -1
$ws.Range("X95").Value2 = -1
$ws.Range("Y95").Value2 = 0.7
$ws.Range("Z95").Value2 = -0.5
$ws.Range("AA95").Value2 = 0.4875
$ws.Range("AB95").Value2 = 0.925
$ws.Range("AC95").Value2 = -1

# Row 96
$ws.Range("B96").Value2 = 6221693
$ws.Range("F96").Value2 = 'Zhetysu'
$ws.Range("G96").Value2 = 'Shakhter Karagandy'
$ws.Range("H96").Value2 = 1
$ws.Range("I96").Value2 = 3
$ws.Range("J96").Value2 = 'A'
$ws.Range("K96").Value2 = 2
$ws.Range("L96").Value2 = 3.4
$ws.Range("M96").Value2 = 3.1
$ws.Range("N96").Value2 = 2.2
$ws.Range("O96").Value2 = 3.3
$ws.Range("P96").Value2 = 2.8
$ws.Range("Q96").Value2 = -0.25
$ws.Range("R96").Value2 = 1.95
$ws.Range("S96").Value2 = 1.85
$ws.Range("T96").Value2 = 2.5
$ws.Range("U96").Value2 = 1.85
$ws.Range("V96").Value2 = 1.95
$ws.Range("W96").Value2 = -1
$ws.Range("X96").Value2 = -1
$ws.Range("Y96").Value2 = 1.8
$ws.Range("Z96").Value2 = -1
$ws.Range("AA96").Value2 = 0.8500000000000001
$ws.Range("AB96").Value2 = 0.8500000000000001
$ws.Range("AC96").Value2 = -1

# Row 104
$ws.Range("B104").Value2 = 6221699
$ws.Range("F104").Value2 = 'FK Maktaaral'
$ws.Range("G104").Value2 = 'Kaisar Kyzylorda'
$ws.Range("H104").Value2 = 2
$ws.Range("I104").Value2 = 2
$ws.Range("J104").Value2 = 'D'
$ws.Range("K104").Value2 = 3.1
$ws.Range("L104").Value2 = 3.2
$ws.Range("M104").Value2 = 2.1
$ws.Range("N104").Value2 = 2.1
$ws.Range("O104").Value2 = 3.1
$ws.Range("P104").Value2 = 3.2
$ws.Range("Q104").Value2 = -0.25
$ws.Range("R104").Value2 = 1.85
$ws.Range("S104").Value2 = 1.95
$ws.Range("T104").Value2 = 2.25
$ws.Range("U104").Value2 = 1.975
$ws.Range("V104").Value2 = 1.725
$ws.Range("W104").Value2 = -1
$ws.Range("X104").Value2 = 2.1
$ws.Range("Y104").Value2 = -1
$ws.Range("Z104").Value2 = -0.5
$ws.Range("AA104").Value2 = 0.475
$ws.Range("AB104").Value2 = 0.9750000000000001
$ws.Range("AC104").Value2 = -1

# Row 105
$ws.Range("B105").Value2 = 6221703
$ws.Range("F105").Value2 = 'Shakhter Karagandy'
$ws.Range("G105").Value2 = 'FK Aktobe'
$ws.Range("H105").Value2 = 0
$ws.Range("I105").Value2 = 1
$ws.Range("J105").Value2 = 'A'
$ws.Range("K105").Value2 = 3.6
$ws.Range("L105").Value2 = 3.5
$ws.Range("M105").Value2 = 1.8
$ws.Range("N105").Value2 = 3.1
$ws.Range("O105").Value2 = 3.5
$ws.Range("P105").Value2 = 1.909
$ws.Range("Q105").Value2 = 0.5
$ws.Range("R105").Value2 = 1.825
$ws.Range("S105").Value2 = 1.975
$ws.Range("T105").Value2 = 2.5
$ws.Range("U105").Value2 = 1.75
$ws.Range("V105").Value2 = 1.95
$ws.Range("W105").Value2 = -1
$ws.Range("X105").Value2 = -1
$ws.Range("Y105").Value2 = 0.909
$ws.Range("Z105").Value2 = -1
$ws.Range("AA105").Value2 = 0.9750000000000001
$ws.Range("AB105").Value2 = -1
$ws.Range("AC105").Value2 = 0.95

# Row 116
$ws.Range("B116").Value2 = 6221712
$ws.Range("F116").Value2 = 'FK Aksu'
$ws.Range("G116").Value2 = 'Shakhter Karagandy'
$ws.Range("H116").Value2 = 2
$ws.Range("I116").Value2 = 1
$ws.Range("J116").Value2 = 'H'
$ws.Range("K116").Value2 = 2.1
$ws.Range("L116").Value2 = 3.25
$ws.Range("M116").Value2 = 3
$ws.Range("N116").Value2 = 2.15
$ws.Range("O116").Value2 = 3.25
$ws.Range("P116").Value2 = 2.9
$ws.Range("Q116").Value2 = -0.25
$ws.Range("R116").Value2 = 1.95
$ws.Range("S116").Value2 = 1.85
$ws.Range("T116").Value2 = 2.5
$ws.Range("U116").Value2 = 1.975
$ws.Range("V116").Value2 = 1.825
$ws.Range("W116").Value2 = 1.15
$ws.Range("X116").Value2 = -1
$ws.Range("Y116").Value2 = -1
$ws.Range("Z116").Value2 = 0.95
$ws.Range("AA116").Value2 = -1
$ws.Range("AB116").Value2 = 0.9750000000000001
$ws.Range("AC116").Value2 = -1

# Row 117
$ws.Range("B117").Value2 = 6221708
$ws.Range("F117").Value2 = 'Kaisar Kyzylorda'
$ws.Range("G117").Value2 = 'Kairat Almaty'
$ws.Range("H117").Value2 = 0
$ws.Range("I117").Value2 = 0
$ws.Range("J117").Value2 = 'D'
$ws.Range("K117").Value2 = 3
$ws.Range("L117").Value2 = 3.4
$ws.Range("M117").Value2 = 2.05
$ws.Range("N117").Value2 = 3.2
$ws.Range("O117").Value2 = 3.4
$ws.Range("P117").Value2 = 1.95
$ws.Range("Q117").Value2 = 0.5
$ws.Range("R117").Value2 = 1.75
$ws.Range("S117").Value2 = 1.95
$ws.Range("T117").Value2 = 2.25
$ws.Range("U117").Value2 = 1.925
$ws.Range("V117").Value2 = 1.875
$ws.Range("W117").Value2 = -1
$ws.Range("X117").Value2 = 2.4
$ws.Range("Y117").Value2 = -1
$ws.Range("Z117").Value2 = 0.75
$ws.Range("AA117").Value2 = -1
$ws.Range("AB117").Value2 = -1
$ws.Range("AC117").Value2 = 0.875

# Row 177
$ws.Range("B177").Value2 = 6221753
$ws.Range("F177").Value2 = 'FK Aksu'
$ws.Range("G177").Value2 = 'Tobol Kostanay'
$ws.Range("H177").Value2 = 0
$ws.Range("I177").Value2 = 3
$ws.Range("J177").Value2 = 'A'
$ws.Range("K177").Value2 = 2.75
$ws.Range("L177").Value2 = 3.1
$ws.Range("M177").Value2 = 2.375
$ws.Range("N177").Value2 = 2.625
$ws.Range("O177").Value2 = 3.2
$ws.Range("P177").Value2 = 2.45
$ws.Range("Q177").Value2 = 0
$ws.Range("R177").Value2 = 2
$ws.Range("S177").Value2 = 1.8
$ws.Range("T177").Value2 = 2.5
$ws.Range("U177").Value2 = 1.9
$ws.Range("V177").Value2 = 1.9
$ws.Range("W177").Value2 = -1
$ws.Range("X177").Value2 = -1
$ws.Range("Y177").Value2 = 1.45
$ws.Range("Z177").Value2 = -1
$ws.Range("AA177").Value2 = 0.8
$ws.Range("AB177").Value2 = 0.8999999999999999
$ws.Range("AC177").Value2 = -1

# Row 178
$ws.Range("B178").Value2 = 6221752
$ws.Range("F178").Value2 = 'FK Kyzylzhar'
$ws.Range("G178").Value2 = 'Kaisar Kyzylorda'
$ws.Range("H178").Value2 = 0
$ws.Range("I178").Value2 = 1
$ws.Range("J178").Value2 = 'A'
$ws.Range("K178").Value2 = 1.833
$ws.Range("L178").Value2 = 3.2
$ws.Range("M178").Value2 = 4
$ws.Range("N178").Value2 = 1.85
$ws.Range("O178").Value2 = 3.2
$ws.Range("P178").Value2 = 4
$ws.Range("Q178").Value2 = -0.5
$ws.Range("R178").Value2 = 1.9
$ws.Range("S178").Value2 = 1.9
$ws.Range("T178").Value2 = 2
$ws.Range("U178").Value2 = 1.775
$ws.Range("V178").Value2 = 2.025
$ws.Range("W178").Value2 = -1
$ws.Range("X178").Value2 = -1
$ws.Range("Y178").Value2 = 3
$ws.Range("Z178").Value2 = -1
$ws.Range("AA178").Value2 = 0.8999999999999999
$ws.Range("AB178").Value2 = -1
$ws.Range("AC178").Value2 = 1.025

# Row 179
$ws.Range("B179").Value2 = 6221815
$ws.Range("F179").Value2 = 'FK Atyrau'
$ws.Range("G179").Value2 = 'Kairat Almaty'
$ws.Range("H179").Value2 = 0
$ws.Range("I179").Value2 = 0
$ws.Range("J179").Value2 = 'D'
$ws.Range("K179").Value2 = 3
$ws.Range("L179").Value2 = 3
$ws.Range("M179").Value2 = 2.25
$ws.Range("N179").Value2 = 3.1
$ws.Range("O179").Value2 = 3.1
$ws.Range("P179").Value2 = 2.15
$ws.Range("Q179").Value2 = 0.25
$ws.Range("R179").Value2 = 1.85
$ws.Range("S179").Value2 = 1.95
$ws.Range("T179").Value2 = 2.25
$ws.Range("U179").Value2 = 1.8
$ws.Range("V179").Value2 = 2
$ws.Range("W179").Value2 = -1
$ws.Range("X179").Value2 = 2.1
$ws.Range("Y179").Value2 = -1
$ws.Range("Z179").Value2 = 0.425
$ws.Range("AA179").Value2 = -0.5
$ws.Range("AB179").Value2 = -1
$ws.Range("AC179").Value2 = 1
